# Change year headers J1: 2020 -> 2022 and K1: 2025 -> 2026
# on the "technological_readiness_bass" sheet (data + code - Bass diffusion
# model columns). This matches the commit:
#   "changed years: 2020, 2025 to 2022, 2026. Both in input files and in code"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("technological_readiness_bass")

# The two year-header cells get left-aligned (Excel mints a new cellXfs
# entry for this the first time it's applied) before the values change.
$ws.Range("J1:K1").HorizontalAlignment = -4131
$ws.Range("J1").Value = 2022
$ws.Range("K1").Value = 2026

# Matches the author's final selection left in the saved file.
$null = $ws.Range("F25").Select()
